$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.186.95'
$ws.Range("E2").Value = '  +2.42%  '

$ws.Range("D3").Value = '2.424.31'
$ws.Range("E3").Value = '  +0.35%  '

$ws.Range("E4").Value = '  +0.47%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.09'
$ws.Range("E5").Value = '  +1.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.21'
$ws.Range("E6").Value = '  +5.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.538'
$ws.Range("E8").Value = '  +0.56%  '

$ws.Range("D9").Value = '2.441.27'
$ws.Range("E9").Value = '  +1.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.110'
$ws.Range("E10").Value = '  +4.33%  '

$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.25'
$ws.Range("E12").Value = '  +3.90%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.350'
$ws.Range("E13").Value = '  +3.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.76'
$ws.Range("E14").Value = '  +4.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000178'
$ws.Range("E15").Value = '  +6.17%  '

$ws.Range("D16").Value = '2.844.69'
$ws.Range("E16").Value = '  +1.20%  '

$ws.Range("D17").Value = '61.826.60'
$ws.Range("E17").Value = '  +1.87%  '

$ws.Range("D18").Value = '2.435.24'
$ws.Range("E18").Value = '  +1.35%  '

$ws.Range("E19").Value = '  -1.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.84'
$ws.Range("E20").Value = '  +2.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.90'
$ws.Range("E21").Value = '  +1.26%  '

$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.04'
$ws.Range("E23").Value = '  +13.64%  '

$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.07'
$ws.Range("E24").Value = '  -1.09%  '

$ws.Range("E25").Value = '  -0.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.47'
$ws.Range("E26").Value = '  +1.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '622.91'
$ws.Range("E27").Value = '  +11.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.45'
$ws.Range("E28").Value = '  +4.00%  '

$ws.Range("D29").Value = '0.0₃0990'
$ws.Range("E29").Value = '  +7.98%  '

$ws.Range("D30").Value = '2.532.11'
$ws.Range("E30").Value = '  +0.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.10'
$ws.Range("E31").Value = '  +2.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.43'
$ws.Range("E32").Value = '  +9.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.83'
$ws.Range("E33").Value = '  +1.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.135'
$ws.Range("E34").Value = '  +3.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("E35").Value = '  +5.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("E36").Value = '  -0.80%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.64'
$ws.Range("E37").Value = '  +2.86%  '

$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.373'
$ws.Range("E38").Value = '  +1.49%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.42'
$ws.Range("E39").Value = '  +6.96%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '151.61'
$ws.Range("E40").Value = '  -0.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.57'
$ws.Range("E41").Value = '  +2.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("E42").Value = '  +17.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.73'
$ws.Range("E43").Value = '  +5.23%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.32'
$ws.Range("E44").Value = '  +2.75%  '

$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.06%  '

$ws.Range("D46").Value = '0.0₆0286'
$ws.Range("E46").Value = '  -1.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '143.65'
$ws.Range("E47").Value = '  +0.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.58'
$ws.Range("E48").Value = '  +2.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.30'
$ws.Range("E49").Value = '  +6.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.600'
$ws.Range("E50").Value = '  +2.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0515'
$ws.Range("E51").Value = '  +3.56%  '
